# Updated symbol list on Tue Dec 27 17:48:42 UTC 2022 with GitHub Actions
#
# The "Price" column (D) in this sheet is stored as TEXT (inline strings
# that merely look numeric), not as real numbers. Assigning a plain
# numeric-looking string via COM .Value makes Excel coerce it to a number,
# so we prefix with an apostrophe to force text entry, then reset the
# cell Style back to "Normal" so the quote-prefix flag doesn't leave a
# stray NumberFormat/style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# Price (column D) updates
Set-TextValue "D2"  "244.74"
Set-TextValue "D3"  "23.94"
Set-TextValue "D4"  "5.361"
Set-TextValue "D5"  "0.05857"
Set-TextValue "D8"  "0.8137"
Set-TextValue "D9"  "0.9236"
Set-TextValue "D10" "0.1409"
Set-TextValue "D11" "0.07403"
Set-TextValue "D12" "0.03088"
Set-TextValue "D13" "0.03059"
Set-TextValue "D14" "0.09365"
Set-TextValue "D15" "3.880"
Set-TextValue "D16" "0.001561"
Set-TextValue "D17" "0.04702"
Set-TextValue "D18" "0.0006049"

# E18 got a "Worstin24h" suffix added
$ws.Range("E18").Value = "17OneONEWorstin24h"

Set-TextValue "D19" "0.005947"
Set-TextValue "D20" "0.001252"
Set-TextValue "D21" "0.004698"
Set-TextValue "D22" "0.00008839"
Set-TextValue "D25" "0.3226"
Set-TextValue "D40" "0.03846"

# Rows 41-43 got reshuffled (KickToken moved up from row43->41,
# BKEXToken moved from row41->42, CEJI moved from row42->43) with
# refreshed prices / labels.
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006458"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1066"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002940"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.008266"
Set-TextValue "D45" "0.00005253"
Set-TextValue "D46" "0.00000000753"
Set-TextValue "D48" "0.001739"
Set-TextValue "D49" "0.00002109"
Set-TextValue "D50" "0.0002008"
